# Append new daily NAV rows (2024-08-28 .. 2024-09-25) to the Pharma_stocks
# tracker sheet, extending the data range from A1:J638 to A1:J659.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 639; Date = "2024-08-28"; C = 2200.75; D = 1539.5; E = 1707.449951171875; F = 1969.050048828125; G = 1138.300048828125; H = 8555.050048828125; I = 0; J = 183.8561105737055 },
    @{ Row = 640; Date = "2024-08-29"; C = 2193.75; D = 1499.150024414062; E = 1691.300048828125; F = 1961.150024414062; G = 1132.050048828125; H = 8477.400146484375; I = -0.009076498898377167; J = 182.1873407886234 },
    @{ Row = 641; Date = "2024-08-30"; C = 2240.199951171875; D = 1537.550048828125; E = 1731.75; F = 1953.800048828125; G = 1127.900024414062; H = 8591.200073242188; I = 0.01342391827581785; J = 184.6330087622584 },
    @{ Row = 642; Date = "2024-09-02"; C = 2232.75; D = 1537.550048828125; E = 1687.900024414062; F = 1970.599975585938; G = 1111.550048828125; H = 8540.35009765625; I = -0.005918844300264037; J = 183.5401947307053 },
    @{ Row = 643; Date = "2024-09-03"; C = 2240.25; D = 1530.599975585938; E = 1687.5; F = 1924.650024414062; G = 1114; H = 8497; I = -0.005075915759957742; J = 182.608560163686 },
    @{ Row = 644; Date = "2024-09-04"; C = 2277.25; D = 1556.550048828125; E = 1686.550048828125; F = 1924.650024414062; G = 1127.900024414062; H = 8572.900146484375; I = 0.008932581674046723; J = 184.2397260417282 },
    @{ Row = 645; Date = "2024-09-05"; C = 2290.199951171875; D = 1555.75; E = 1709.449951171875; F = 1933.599975585938; G = 1115.150024414062; H = 8604.14990234375; I = 0.00364517903223101; J = 184.9113128279995 },
    @{ Row = 646; Date = "2024-09-06"; C = 2256.5; D = 1559.900024414062; E = 1702.699951171875; F = 1928.400024414062; G = 1100; H = 8547.5; I = -0.006584020848860234; J = 183.6938528891499 },
    @{ Row = 647; Date = "2024-09-09"; C = 2216.800048828125; D = 1546.25; E = 1704.199951171875; F = 1937.099975585938; G = 1104.150024414062; H = 8508.5; I = -0.004562737642585551; J = 182.855706031861 },
    @{ Row = 648; Date = "2024-09-10"; C = 2222.550048828125; D = 1545.550048828125; E = 1727.849975585938; F = 1912.150024414062; G = 1113.199951171875; H = 8521.300048828125; I = 0.001504383713712758; J = 183.1307911779747 },
    @{ Row = 649; Date = "2024-09-11"; C = 2209.39990234375; D = 1591.949951171875; E = 1725.650024414062; F = 1867.75; G = 1112.599975585938; H = 8507.349853515625; I = -0.001637097066476197; J = 182.8309882969558 },
    @{ Row = 650; Date = "2024-09-12"; C = 2247.5; D = 1592.849975585938; E = 1747.949951171875; F = 1883.349975585938; G = 1120.099975585938; H = 8591.749877929688; I = 0.009920836202496664; J = 184.6448245845905 },
    @{ Row = 651; Date = "2024-09-13"; C = 2256.449951171875; D = 1582.5; E = 1753.699951171875; F = 1923.300048828125; G = 1118.550048828125; H = 8634.5; I = 0.004975717714982386; J = 185.5635651092558 },
    @{ Row = 652; Date = "2024-09-16"; C = 2251.85009765625; D = 1577.75; E = 1741.449951171875; F = 1900.949951171875; G = 1115.849975585938; H = 8587.849975585938; I = -0.005402747630327465; J = 184.5610119975867 },
    @{ Row = 653; Date = "2024-09-17"; C = 2270.39990234375; D = 1561.699951171875; E = 1713; F = 1875.599975585938; G = 1110.949951171875; H = 8531.649780273438; I = -0.006544151967287428; J = 183.3532166878381 },
    @{ Row = 654; Date = "2024-09-18"; C = 2224.949951171875; D = 1543.050048828125; E = 1646.050048828125; F = 1857; G = 1079.949951171875; H = 8351; I = -0.02117407358786916; J = 179.4708821851173 },
    @{ Row = 655; Date = "2024-09-19"; C = 2171.89990234375; D = 1515.050048828125; E = 1649.800048828125; F = 1886.5; G = 1054.449951171875; H = 8277.699951171875; I = -0.00877739777608969; J = 177.8955948629528 },
    @{ Row = 656; Date = "2024-09-20"; C = 2151.699951171875; D = 1481.099975585938; E = 1636.75; F = 1897.25; G = 1054.599975585938; H = 8221.39990234375; I = -0.0068014121265841; J = 176.685653606786 },
    @{ Row = 657; Date = "2024-09-23"; C = 2182.25; D = 1440.400024414062; E = 1712.449951171875; F = 1952; G = 1055.25; H = 8342.349975585938; I = 0.01471161537923817; J = 179.2849849856784 },
    @{ Row = 658; Date = "2024-09-24"; C = 2215.75; D = 1414.25; E = 1697.5; F = 1944.349975585938; G = 1051.550048828125; H = 8323.400024414062; I = -0.002271536344954651; J = 178.8777326261788 },
    @{ Row = 659; Date = "2024-09-25"; C = 2221.10009765625; D = 1416.400024414062; E = 1689.199951171875; F = 1909.550048828125; G = 1063.449951171875; H = 8299.700073242188; I = -0.002847388219040138; J = 178.3683982776504 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Column A holds the trade date as plain text (matches the existing
    # rows above it), not an Excel date serial - force text with a leading
    # apostrophe, then drop the resulting "Text" number-format style so the
    # cell stays styleless like its neighbours.
    $ws.Range("A$row").Value = "'" + $r.Date
    $ws.Range("A$row").ClearFormats()

    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
}
